$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.782.93'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '2.129.02'
$ws.Range("E3").Value = '  +10.64%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'255.92"
$ws.Range("E5").Value = '  +2.00%  '
$ws.Range("D6").Value = "'0.668"
$ws.Range("E6").Value = '  -4.28%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'47.63"
$ws.Range("E8").Value = '  +7.29%  '
$ws.Range("D9").Value = "'60.71"
$ws.Range("E9").Value = '  +2.92%  '
$ws.Range("D10").Value = "'0.375"
$ws.Range("E10").Value = '  +1.93%  '
$ws.Range("D11").Value = "'0.0743"
$ws.Range("E11").Value = '  -3.52%  '
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("D13").Value = '2.433.11'
$ws.Range("E13").Value = '  +10.35%  '
$ws.Range("D14").Value = "'14.49"
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").Value = "'0.847"
$ws.Range("E15").Value = '  +5.67%  '
$ws.Range("D16").Value = '2.124.34'
$ws.Range("E16").Value = '  +10.36%  '
$ws.Range("D17").Value = "'5.15"
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '36.796.06'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").Value = "'73.95"
$ws.Range("E19").Value = '  -0.82%  '
$ws.Range("D20").Value = '0.0₃0843'
$ws.Range("E20").Value = '  -3.05%  '
$ws.Range("D21").Value = "'13.36"
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").Value = "'242.27"
$ws.Range("E22").Value = '  -3.72%  '
$ws.Range("D23").Value = "'5.23"
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = '  -7.44%  '
$ws.Range("D26").Value = "'172.00"
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("D27").Value = "'21.70"
$ws.Range("E27").Value = '  +15.18%  '
$ws.Range("D28").Value = "'9.29"
$ws.Range("E28").Value = '  +5.46%  '
$ws.Range("D29").Value = "'2.04"
$ws.Range("E29").Value = '  -8.42%  '
$ws.Range("D30").Value = "'28.32"
$ws.Range("E30").Value = '  +56.39%  '
$ws.Range("D31").Value = "'0.124"
$ws.Range("E31").Value = '  -4.42%  '
$ws.Range("D32").Value = "'4.53"
$ws.Range("E32").Value = '  -0.27%  '
$ws.Range("D33").Value = "'0.0952"
$ws.Range("E33").Value = '  +11.92%  '
$ws.Range("D34").Value = "'0.0602"
$ws.Range("E34").Value = '  -1.42%  '
$ws.Range("D35").Value = "'2.38"
$ws.Range("E35").Value = '  +17.66%  '
$ws.Range("E36").Value = '  -4.67%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = "'0.948"
$ws.Range("E37").Value = '  +9.43%  '
$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").Value = "'4.19"
$ws.Range("E39").Value = '  -3.67%  '
$ws.Range("D41").Value = "'1.20"
$ws.Range("E41").Value = '  +7.88%  '
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("D43").Value = "'99.34"
$ws.Range("E43").Value = '  -7.81%  '
$ws.Range("E44").Value = '  +12.05%  '
$ws.Range("D45").Value = "'16.26"
$ws.Range("E45").Value = '  -5.39%  '
$ws.Range("D46").Value = '1.360.46'
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").Value = "'0.0845"
$ws.Range("E47").Value = '  +4.03%  '
$ws.Range("D48").Value = "'7.18"
$ws.Range("E48").Value = '  +11.18%  '
$ws.Range("D49").Value = '2.314.67'
$ws.Range("E49").Value = '  +9.74%  '
$ws.Range("E50").Value = '  -2.84%  '
$ws.Range("D51").Value = "'2.84"
$ws.Range("E51").Value = '  +1.07%  '

# Cells above got an Excel "quote prefix" to keep numeric-looking text
# as text instead of being auto-converted to a number; clear the
# resulting style tweak so formatting stays untouched.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
